$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the row above for the styled columns (A and G)
# before putting the new values in, so the new cells reuse the existing
# cellXf (s="1") instead of minting a new one.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)

$ws.Range("A9").Value = 42654.746493055558
$ws.Range("B9").Value = $false
$ws.Range("C9").Value = 9894.34
$ws.Range("D9").Value = 9909.2000000000007
$ws.Range("E9").Value = 104.43
$ws.Range("F9").Value = 104.74
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = 0.3
$ws.Range("I9").Value = $false

$ws.Columns.Item(1).ColumnWidth = 15.375
